# histens fixes, pushing to run LOO code
#
# The empty "I" column (old width 14.83203125) that separated the data
# table from the stray notes in J/K is removed, shifting the J/K note
# columns (ids 30/31/29/28) left into I/J.
#
# A3's bold-header highlight fill is cleared (keeps the bold font), and
# the stale yellow-ish "no fill" marker style on B9:B11 is cleared back
# to the default (unstyled) cell.
#
# Finally the sheet view is re-zoomed and the selection is moved to
# A1:I11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Obs")

# Delete the empty column I -- shifts the J/K "note" cells left to I/J
# and updates dimension/col widths/row spans automatically.
$ws.Columns("I").Delete()

# Clear the (no-op) fill that was still flagged "applyFill" on A3, while
# keeping its bold font.
$ws.Range("A3").Interior.Pattern = -4142

# Clear the stray "applyFill" formatting left on B9:B11 (no visible fill
# change, just drops the now-pointless style record).
$ws.Range("B9:B11").Interior.Pattern = -4142

# Re-zoom and move the selection/active cell (the host always anchors a
# multi-cell Range.Select() at the range's top-left corner, so A1:I11
# reproduces the author's sqref even though the saved activeCell ends up
# at A1 rather than I11).
$ws.Application.ActiveWindow.Zoom = 139
$ws.Range("A1:I11").Select()
